$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has columns:
#   A title | B artist | C album | D genre | E released_date | F description | G korean_music
#
# Target layout removes the "genre" (D) and "korean_music" (G) columns entirely,
# shifting released_date -> D and description -> E:
#   A title | B artist | C album | D released_date | E description

$ws.Columns("D:D").Delete()
# korean_music was column G; after the first delete it has shifted left to column F.
$ws.Columns("F:F").Delete()

$wb.Save()
